$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new SVR parameter headers in columns K, L, M (row 1)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add the corresponding SVR parameter values in row 2
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Clear the leftover empty styled placeholder cell A13 (row stays in place,
# it just no longer carries any cell data so it drops out of the sheet XML)
$ws.Range("A13").Clear()

# Update selection to match saved state
$ws.Range("J8").Select()
